$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1597
$ws.Range("F5").Value = 9135
$ws.Range("F6").Value = 272
$ws.Range("F7").Value = 117
$ws.Range("F8").Value = 1275
$ws.Range("F10").Value = 636
$ws.Range("F13").Value = 151
$ws.Range("F17").Value = 1510
$ws.Range("F18").Value = 1328
$ws.Range("F22").Value = 86
$ws.Range("F23").Value = 237
$ws.Range("F26").Value = 67
$ws.Range("F28").Value = 315
$ws.Range("F29").Value = 315
$ws.Range("F30").Value = 1073
$ws.Range("F34").Value = 203
$ws.Range("E36").Value = "2024.06.22 10:00-06.23 16:20"
$ws.Range("F39").Value = 136
$ws.Range("F41").Value = 155
$ws.Range("F42").Value = 67
$ws.Range("F47").Value = 47

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F25").Value = 21
$ws.Range("F27").Value = 234
$ws.Range("F30").Value = 231
$ws.Range("F31").Value = 4
$ws.Range("F32").Value = 153

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 315
$ws.Range("F6").Value = 146
$ws.Range("F7").Value = 2109
$ws.Range("F8").Value = 3167

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1597
$ws.Range("F6").Value = 9135
$ws.Range("F7").Value = 315
$ws.Range("F8").Value = 146
$ws.Range("F9").Value = 272
$ws.Range("F10").Value = 2109
$ws.Range("F11").Value = 3167
$ws.Range("F12").Value = 117
$ws.Range("F13").Value = 1275
$ws.Range("F14").Value = 636
$ws.Range("F17").Value = 151
$ws.Range("F19").Value = 1510
$ws.Range("F21").Value = 1328
$ws.Range("F24").Value = 237
$ws.Range("F26").Value = 67
$ws.Range("F27").Value = 315
$ws.Range("F28").Value = 315
$ws.Range("F29").Value = 1073
$ws.Range("F35").Value = 203
$ws.Range("F37").Value = 234
$ws.Range("E38").Value = "2024.06.22 10:00-06.23 16:20"
$ws.Range("F41").Value = 136
$ws.Range("F42").Value = 231
$ws.Range("F43").Value = 155
$ws.Range("F44").Value = 67
$ws.Range("F45").Value = 153
$ws.Range("F50").Value = 47
